$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1169
$ws.Cells.Item(3, 6).Value = 1981
$ws.Cells.Item(4, 6).Value = 627
$ws.Cells.Item(5, 6).Value = 1279
$ws.Cells.Item(6, 6).Value = 72
$ws.Cells.Item(7, 6).Value = 53
$ws.Cells.Item(9, 6).Value = 348
$ws.Cells.Item(10, 6).Value = 132
$ws.Cells.Item(11, 6).Value = 106
$ws.Cells.Item(12, 6).Value = 869
$ws.Cells.Item(13, 6).Value = 269
$ws.Cells.Item(14, 6).Value = 138
$ws.Cells.Item(17, 6).Value = 348
$ws.Cells.Item(18, 6).Value = 261
$ws.Cells.Item(21, 6).Value = 676
$ws.Cells.Item(22, 6).Value = 209
$ws.Cells.Item(24, 6).Value = 921
$ws.Cells.Item(25, 6).Value = 381
$ws.Cells.Item(26, 6).Value = 202
$ws.Cells.Item(28, 6).Value = 313

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 338
$ws.Cells.Item(7, 6).Value = 266
$ws.Cells.Item(9, 6).Value = 5
$ws.Cells.Item(11, 6).Value = 132
$ws.Cells.Item(12, 6).Value = 28

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 333

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 333
$ws.Cells.Item(3, 6).Value = 1169
$ws.Cells.Item(4, 6).Value = 1981
$ws.Cells.Item(5, 6).Value = 627
$ws.Cells.Item(6, 6).Value = 1279
$ws.Cells.Item(7, 6).Value = 72
$ws.Cells.Item(9, 6).Value = 53
$ws.Cells.Item(11, 6).Value = 348
$ws.Cells.Item(12, 6).Value = 132
$ws.Cells.Item(13, 6).Value = 106
$ws.Cells.Item(14, 6).Value = 869
$ws.Cells.Item(15, 6).Value = 269
$ws.Cells.Item(16, 6).Value = 138
$ws.Cells.Item(19, 6).Value = 338
$ws.Cells.Item(22, 6).Value = 348
$ws.Cells.Item(24, 6).Value = 266
$ws.Cells.Item(25, 6).Value = 261
$ws.Cells.Item(28, 6).Value = 676
$ws.Cells.Item(29, 6).Value = 209
$ws.Cells.Item(31, 6).Value = 921
$ws.Cells.Item(32, 6).Value = 381
$ws.Cells.Item(34, 6).Value = 5
$ws.Cells.Item(35, 6).Value = 202
$ws.Cells.Item(37, 6).Value = 313
$ws.Cells.Item(39, 6).Value = 132
$ws.Cells.Item(42, 6).Value = 28
